$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow the data edits below.
$ws.Unprotect()

# --- Update the "as of" date in the confidential disclosure footer (A37) ---
$footerText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."
$ws.Range("A37").Value = $footerText
# Re-applying the original (non-custom) row height avoids an incidental
# row-height side effect from the embedded line break.
$ws.Rows.Item(37).AutoFit()

# --- Refresh Weight (D) and Percent Change (E) figures for each holding ---
$ws.Range("D2").Value = 0.03867380301545728
$ws.Range("E2").Value = 0.003123779773525959
$ws.Range("D3").Value = 0.02189576162602871
$ws.Range("E3").Value = 0.002018779342723098
$ws.Range("D4").Value = 0.02016578471331853
$ws.Range("E4").Value = -0.0004139072847683112
$ws.Range("D5").Value = 0.04074758977613847
$ws.Range("E5").Value = 0.01016117729502453
$ws.Range("D6").Value = 0.03771183832662866
$ws.Range("E6").Value = 0.0003916960438699046
$ws.Range("D7").Value = 0.02113186787181619
$ws.Range("E7").Value = 0.001948937828883057
$ws.Range("D8").Value = 0.03767191662730347
$ws.Range("E8").Value = -0.0004518752824218542
$ws.Range("D9").Value = 0.02157550942456459
$ws.Range("E9").Value = 0.002199010445299576
$ws.Range("D10").Value = 0.02581456788550845
$ws.Range("E10").Value = 0.003935332907891942
$ws.Range("D11").Value = 0.0239351728932057
$ws.Range("E11").Value = -0.001720676799541132
$ws.Range("D12").Value = 0.05826506944688483
$ws.Range("E12").Value = -0.00627838800186209
$ws.Range("D13").Value = 0.02660860883769832
$ws.Range("E13").Value = 0.003343239227340256
$ws.Range("D14").Value = 0.02730084121966964
$ws.Range("E14").Value = 0.0009896091044037103
$ws.Range("D15").Value = 0.03500106159047012
$ws.Range("E15").Value = 0.001757160428747184
$ws.Range("D16").Value = 0.01883469531628657
$ws.Range("E16").Value = 0.00707013574660631
$ws.Range("D17").Value = 0.02990283544363673
$ws.Range("E17").Value = 0.006427325314479759
$ws.Range("D18").Value = 0.02418002964326203
$ws.Range("E18").Value = 0.003242992819087354
$ws.Range("D19").Value = 0.1341685395796453
$ws.Range("E19").Value = 0.004032258064516014
$ws.Range("D20").Value = 0.00963288090415609
$ws.Range("E20").Value = -0.0003933394519469724
$ws.Range("D21").Value = 0.0159160731441641
$ws.Range("E21").Value = 0.001231705549920248
$ws.Range("D22").Value = 0.01724677815069361
$ws.Range("E22").Value = 0.0005412719891744189
$ws.Range("D23").Value = 0.01711855646166434
$ws.Range("E23").Value = -0.01398601398601407
$ws.Range("D24").Value = 0.02145201024742248
$ws.Range("E24").Value = -0.004623004953219434
$ws.Range("D25").Value = 0.0116372028097529
$ws.Range("E25").Value = 0.006006955422067639
$ws.Range("D26").Value = 0.04340702292375637
$ws.Range("E26").Value = -0.005408174883297145
$ws.Range("D27").Value = 0.02574999028109799
$ws.Range("E27").Value = 0.00009809691975681822
$ws.Range("D28").Value = 0.04792037146465793
$ws.Range("E28").Value = 0.004637539663168155
$ws.Range("D29").Value = 0.05726957217801092
$ws.Range("E29").Value = -0.00864499154294307
$ws.Range("D30").Value = 0.01311727098294607
$ws.Range("E30").Value = 0.01806810284920068
$ws.Range("D31").Value = 0.01464802378953282
$ws.Range("E31").Value = -0.003280224929709474
$ws.Range("D32").Value = 0.04455162801417445
$ws.Range("E32").Value = 0.001042752867570274
$ws.Range("D33").Value = 0.01674712541044635
$ws.Range("E33").Value = -0.001767352185089943
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 0.0009311206457833787

# Re-apply sheet protection (structure only; original password cannot be
# round-tripped through the object model, so we restore the protected
# state without attempting to recreate the exact legacy hash).
$ws.Protect()
